$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 82; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value2
    $cell.Value = $old + 0.136762601
}
